$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-coerced to numbers by Excel (the Price column is text-typed).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "56.407.13"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "2.314.38"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "515.75"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "131.65"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  -0.96%  "
$ws.Range("E9").Value = "  -2.53%  "
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").Value = "0.336"
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "23.45"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.727.96"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "56.355.47"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("D17").Value = "2.317.18"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").Value = "10.38"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").Value = "329.54"
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("D20").Value = "4.14"
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("D21").Value = "6.69"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "60.91"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").Value = "8.63"
$ws.Range("E25").Value = "  +8.69%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("E27").Value = "  +1.74%  "
$ws.Range("D28").Value = "168.31"
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("D29").Value = "1.68"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("D30").Value = "0.0₃0717"
$ws.Range("E30").Value = "  -2.73%  "
$ws.Range("D31").Value = "6.13"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").Value = "18.29"
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("D36").Value = "3.92"
$ws.Range("E36").Value = "  -2.51%  "
$ws.Range("D37").Value = "0.880"
$ws.Range("E37").Value = "  -4.29%  "
$ws.Range("D38").Value = "1.57"
$ws.Range("E38").Value = "  +1.66%  "
$ws.Range("D39").Value = "38.63"
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("D40").Value = "148.48"
$ws.Range("E40").Value = "  +6.45%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "0.373"
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "3.58"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "284.12"
$ws.Range("E43").Value = "  +1.55%  "
$ws.Range("D44").Value = "5.01"
$ws.Range("E44").Value = "  -4.20%  "
$ws.Range("D45").Value = "0.0926"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").Value = "0.0494"
$ws.Range("E46").Value = "  -2.20%  "
$ws.Range("D47").Value = "0.553"
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("D48").Value = "18.09"
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("B50").Value = "Polygon"
$ws.Range("C50").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D50").Value = "0.375"
$ws.Range("E50").Value = "  -1.38%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "17.07"
$ws.Range("E51").Value = "  +1.60%  "
